$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.604.36'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.05%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.916.55'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +2.65%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.34'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.54%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.07%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5161'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +3.40%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3999'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.42%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.09812'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -2.30%  '

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +2.47%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.27'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +2.12%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.515'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.22%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.25'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.03%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.920.69'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.93%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.479'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.10%  '

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.09%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '94.80'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.22%  '

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.71%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06658'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.10%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.29'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +5.10%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.02%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.311'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +3.53%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.646.55'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.97%  '

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.21%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.313'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.17%  '

$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.390'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -1.43%  '

$ws.Range('B27').NumberFormat = "@"
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').NumberFormat = "@"
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.684'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +6.31%  '

$ws.Range('B28').NumberFormat = "@"
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').NumberFormat = "@"
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.141.29'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +3.01%  '

$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '21.30'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.36%  '

$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '157.96'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.07%  '

$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '129.26'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.29%  '

$ws.Range('B32').NumberFormat = "@"
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').NumberFormat = "@"
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.116'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +5.57%  '

$ws.Range('B33').NumberFormat = "@"
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').NumberFormat = "@"
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.1074'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.42%  '

$ws.Range('B34').NumberFormat = "@"
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').NumberFormat = "@"
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.731'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.59%  '

$ws.Range('B35').NumberFormat = "@"
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').NumberFormat = "@"
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.631'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.88%  '

$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '9.910'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +6.00%  '

$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06790'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.47%  '

$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02441'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +1.84%  '

$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.273'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +5.36%  '

$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2234'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +2.77%  '

$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '11.84'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +2.72%  '

$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.6508'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +3.10%  '

$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'InternetComputer(DFINITY)'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.090'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.31%  '

$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.187'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.74%  '

$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.09%  '

$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '13.72'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.66%  '

$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.6110'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.61%  '

$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.787'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +3.03%  '

$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'WEMIXTOKEN'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.287'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.30%  '

$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.070'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +3.99%  '

$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '125.09'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.04%  '

